$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header cells (P1, Q1), copying the formatting (border/bold/centered)
# from the existing header cell O1 so they share the same style index.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25: swap values in columns I/K and M/O, and populate the
# two new columns P and Q.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1 (was 2)
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2 (new)
}
